$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New")
$ws.Activate()

# --- Fix the "Edge" column (column C) ---
# Rows 2-37 currently hold "Yes" and should become "No".
# Rows 38-73 currently hold "No" and should become "Yes".
$ws.Range("C2:C37").Value = "No"
$ws.Range("C38:C73").Value = "Yes"

# --- Update selection / view state ---
[void]$ws.Range("C38:C73").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1

# --- Re-apply the AutoFilter so its reference covers the full data range ---
$r = $ws.Range("A1:V73")
[void]$r.AutoFilter()
[void]$r.AutoFilter()

# --- Update the hidden _FilterDatabase defined name to match the new range ---
$names = $ws.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "New!_FilterDatabase") {
        $n.RefersTo = "=New!`$A`$1:`$V`$73"
    }
}
